$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: cardholder name / account number / surname
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 20.12.2024"

# Row 6
$ws.Range("B6").Value = "23.12."
$ws.Range("C6").Value = "24.12."
$ws.Range("D6").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E6").Value = "24,95-"

# Row 7
$ws.Range("B7").Value = "27.12."
$ws.Range("C7").Value = "28.12."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-24932818"
$ws.Range("E7").Value = "54,68-"

# Row 8
$ws.Range("B8").Value = "28.12."
$ws.Range("C8").Value = "29.12."
$ws.Range("D8").Value = "ZALANDO MKTPLC EU WFAZII"
$ws.Range("E8").Value = "90,60-"

# Row 9
$ws.Range("B9").Value = "31.12."
$ws.Range("C9").Value = "01.01."
$ws.Range("D9").Value = "KARTENZ./31.12 EDEKA RO"
$ws.Range("E9").Value = "76,60-"

# Row 10
$ws.Range("B10").Value = "03.01."
$ws.Range("C10").Value = "04.01."
$ws.Range("D10").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E10").Value = "78,92-"

# Row 11: cleared out (was a transaction row, now blank)
$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("E11").Style = "Normal"

# Row 12: closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 05.01.2025"
$ws.Range("E12").Value = "325,75-"

# Row 13: next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 14.01.2025"
